# Applies the edits described by the commit "google analytics and update
# report project":
#   1. Remove the "-0.2" penalty text from the "Đăng kí thành giáo viên"
#      row (table 1, row 26, column 3) leaving the paragraph empty.
#   2. Trim the "Sử dụng Google Analytics/Sentry/Crashlytics…" cell text
#      down to just "...Analytics" (table 1, row 60, column 1).
#   3. Fill in the previously empty bonus-score cell next to it with the
#      text "0,25" typed as three separate runs ("0,", "2", "5") that all
#      carry <w:lang w:val="en-US"/>, matching how Word records text
#      typed in several keystrokes/sessions.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Edit 1: delete the "-0.2" run, table row 26 / column 3 -----------
$cell = $t.Cell(26, 3)
$cellRange = $cell.Range
$delRange = $d.Range($cellRange.Start, $cellRange.End - 1)
$delRange.Delete()

# --- Edit 2: trim "/Sentry/Crashlytics…" off the Analytics cell -------
$t = $d.Tables.Item(1)
$cell = $t.Cell(60, 1)
$cellRange = $cell.Range
$cellEnd = $cellRange.End - 1
$findRange = $cell.Range
$findRange.Find.Execute("Analytics", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "", 0)
$tailRange = $d.Range($findRange.End, $cellEnd)
$tailRange.Delete()

# --- Edit 3: type "0,25" (as three runs) into the now-empty score cell
$t = $d.Tables.Item(1)
$cell = $t.Cell(60, 3)
$cellRange = $cell.Range

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>' + `
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' + `
  '<w:body><w:p w14:paraId="521343E1" w14:textId="77777777" w:rsidR="00D66A42" w:rsidRDefault="00D66A42" w:rsidP="00D66A42">' + `
  '<w:pPr><w:widowControl w:val="0"/><w:spacing w:before="144" w:after="144"/><w:jc w:val="center"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' + `
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>0,</w:t></w:r>' + `
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>2</w:t></w:r>' + `
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>5</w:t></w:r>' + `
  '</w:p></w:body></w:document>' + `
  '</pkg:xmlData></pkg:part></pkg:package>'

$cellRange.InsertXML($xml)
